$d = $word.ActiveDocument

# --- Edit 1: merge the "ost" + "):" runs into a single "ost):" run ---
# Locate "Post):" (the end of "Blog Post 3 (Project 1 Closing Blog Post):")
# and narrow the range to just "ost):" (skip the leading "P", which is its
# own separate run that must stay untouched).
$full = $d.Content
$found = $full.Find.Execute("Post):")
if ($found) {
    $merge = $d.Range($full.Start + 1, $full.End)

    # capture the existing (identical) run formatting before we touch it
    $fontName   = $merge.Font.Name
    $fontNameBi = $merge.Font.NameBi
    $bold       = $merge.Font.Bold
    $boldBi     = $merge.Font.BoldBi
    $color      = $merge.Font.Color

    # delete + reinsert forces Word to collapse the two source runs into one
    $merge.Delete()
    $merge.InsertAfter("ost):")

    # re-apply the captured formatting so the new run matches the original
    $merge.Font.Name   = $fontName
    $merge.Font.NameBi = $fontNameBi
    $merge.Font.Bold   = $bold
    $merge.Font.BoldBi = $boldBi
    $merge.Font.Color  = $color
}

# --- Edit 2: drop the block of empty paragraphs trailing the final "Michael" ---
# Find the last paragraph that is exactly "Michael" (the sign-off before the
# run of blank paragraphs at the very end of the document) and remove all of
# the blank paragraphs that follow it except for the very last one.
$count = $d.Paragraphs.Count
$michaelIndex = 0
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Michael`r") {
        $michaelIndex = $i
    }
}

if ($michaelIndex -gt 0) {
    $lastIndex = $d.Paragraphs.Count
    # keep the paragraph right after "Michael" deleted... actually keep the
    # very last paragraph of the document; remove everything blank in-between.
    $firstBlank = $michaelIndex + 1
    $lastBlankToRemove = $lastIndex - 1
    if ($lastBlankToRemove -ge $firstBlank) {
        $p1 = $d.Paragraphs.Item($firstBlank)
        $p2 = $d.Paragraphs.Item($lastBlankToRemove)
        $delRange = $d.Range($p1.Range.Start, $p2.Range.End)
        $delRange.Delete()
    }
}
